$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New POC entries appended to the list (column A), mirroring rows already
# present for "integration testing" etc. higher up the sheet.
$values = @(
    "file upload",
    "spring boot dev tools",
    "injecting custom properties",
    "custom properties with configuraton file",
    "spring boot 2 activemq"
)

$startRow = 39
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Match the scrolled/selected view state saved with the edit.
$excel.Goto($ws.Range("A17"), $true)
$ws.Range("A44").Select()
